# NBA2324.xlsx — append the 2024-02-10 slate of games (11 games, row 778-788)
# to the bottom of the existing results table on Sheet1.
#
# Column layout (row 1 header): A=Away team, B=Away Pts, C=Home team,
# D=Home Pts, E=Overtime, F=Attend., G=Arena, H=Win, I=Loss

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Away team, Away Pts, Home team, Home Pts, Overtime, Attend., Arena, Win, Loss
$newGames = @(
    @("Oklahoma City Thunder", 111, "Dallas Mavericks",       146, "No", 17832, "American Airlines Center",        "Dallas Mavericks",     "Oklahoma City Thunder"),
    @("Detroit Pistons",       106, "Los Angeles Clippers",   112, "No", 17832, "Crypto.com Arena",                "Los Angeles Clippers", "Detroit Pistons"),
    @("San Antonio Spurs",     103, "Brooklyn Nets",          123, "No", 17832, "Barclays Center",                 "Brooklyn Nets",        "San Antonio Spurs"),
    @("Memphis Grizzlies",     106, "Charlotte Hornets",      115, "No", 17832, "Spectrum Center",                 "Charlotte Hornets",    "Memphis Grizzlies"),
    @("Chicago Bulls",         108, "Orlando Magic",          114, "No", 17832, "Amway Center",                    "Orlando Magic",        "Chicago Bulls"),
    @("Philadelphia 76ers",    119, "Washington Wizards",     113, "No", 17832, "Capital One Arena",               "Philadelphia 76ers",   "Washington Wizards"),
    @("Houston Rockets",       113, "Atlanta Hawks",          122, "No", 17832, "State Farm Arena",                "Atlanta Hawks",        "Houston Rockets"),
    @("Indiana Pacers",        125, "New York Knicks",        111, "No", 17832, "Madison Square Garden (IV)",      "Indiana Pacers",       "New York Knicks"),
    @("Cleveland Cavaliers",   119, "Toronto Raptors",         95, "No", 17832, "Scotiabank Arena",                "Cleveland Cavaliers",  "Toronto Raptors"),
    @("Phoenix Suns",          112, "Golden State Warriors",  113, "No", 17832, "Chase Center",                    "Golden State Warriors","Phoenix Suns"),
    @("New Orleans Pelicans",   93, "Portland Trail Blazers",  84, "No", 17832, "Moda Center",                     "New Orleans Pelicans", "Portland Trail Blazers")
)

$firstNewRow = 778
$row = $firstNewRow
foreach ($game in $newGames) {
    $ws.Cells.Item($row, 1).Value = $game[0]

    $ws.Cells.Item($row, 2).Value = $game[1]
    $ws.Cells.Item($row, 2).NumberFormat = "#,##0"

    $ws.Cells.Item($row, 3).Value = $game[2]

    $ws.Cells.Item($row, 4).Value = $game[3]
    $ws.Cells.Item($row, 4).NumberFormat = "#,##0"

    $ws.Cells.Item($row, 5).Value = $game[4]
    $ws.Cells.Item($row, 6).Value = $game[5]
    $ws.Cells.Item($row, 7).Value = $game[6]
    $ws.Cells.Item($row, 8).Value = $game[7]
    $ws.Cells.Item($row, 9).Value = $game[8]

    $row = $row + 1
}

$lastNewRow = $row - 1

# Leave the selection/scroll position where the author left it after typing
# in the new rows (bottom of the sheet, cell A778 active).
try {
    $ws.Range("A" + $firstNewRow).Select()
    $excel.ActiveWindow.ScrollRow = 758
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Cosmetic view-state only; ignore if unsupported by the host.
}
